$wb = $excel.ActiveWorkbook
$styles = $wb.Styles
$good = $styles.Add("Good")
$good.Font.Color = 0x006100
$good.Interior.Color = 0xC6EFCE
$good.Interior.PatternType = -4124
Write-Host "created"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$new.Name = "Average percentuale variation"
$new.Range("B11").Value = 0.5
$new.Range("B11").Style = "Good"
$new.Range("B11").NumberFormat = "0.00%"
